$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Temporarily mark the target cells as Text so that values such as
# "294.00", "0.001270" or "-4.55%" are stored as literal strings rather than
# being auto-converted by Excel into numbers/percentages (which would drop
# trailing zeros or change the stored representation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Step 2: Write the new values (coin names, links, prices and volume%).
$ws.Range("D2").Value = "294.00"
$ws.Range("E2").Value = "-4.55%"
$ws.Range("D3").Value = "40.28"
$ws.Range("E3").Value = "-1.89%"
$ws.Range("D4").Value = "5.065"
$ws.Range("E4").Value = "-2.72%"
$ws.Range("D5").Value = "0.07457"
$ws.Range("E5").Value = "-2.82%"
$ws.Range("D6").Value = "1.596"
$ws.Range("E6").Value = "-1.95%"
$ws.Range("D7").Value = "0.9363"
$ws.Range("E7").Value = "2.38%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "0.1192"
$ws.Range("E8").Value = "-1.90%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1763"
$ws.Range("E9").Value = "-3.15%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.08836"
$ws.Range("E10").Value = "-2.93%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.04182"
$ws.Range("E11").Value = "-0.78%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.1055"
$ws.Range("E12").Value = "0.29%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001270"
$ws.Range("E13").Value = "0.89%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "0.005877"
$ws.Range("E14").Value = "2.66%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "3.376"
$ws.Range("E15").Value = "1.07%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "4.341"
$ws.Range("E16").Value = "0.88%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "2.401"
$ws.Range("E17").Value = "-1.76%"
$ws.Range("D18").Value = "0.3351"
$ws.Range("E18").Value = "0.46%"
$ws.Range("D19").Value = "7.680"
$ws.Range("E19").Value = "4.59%"
$ws.Range("D20").Value = "0.1360"
$ws.Range("E20").Value = "-1.66%"
$ws.Range("D21").Value = "0.2822"
$ws.Range("E21").Value = "4.00%"
$ws.Range("D22").Value = "0.03861"
$ws.Range("E22").Value = "-3.85%"
$ws.Range("D23").Value = "0.001295"
$ws.Range("E23").Value = "2.51%"
$ws.Range("D24").Value = "0.003528"
$ws.Range("E24").Value = "-17.21%"
$ws.Range("D25").Value = "0.0001307"
$ws.Range("E25").Value = "0.45%"
$ws.Range("D26").Value = "0.0003748"
$ws.Range("E26").Value = "-95.01%"
$ws.Range("D38").Value = "0.02334"
$ws.Range("E38").Value = "-6.59%"
$ws.Range("D39").Value = "0.05054"
$ws.Range("E39").Value = "-4.96%"
$ws.Range("D40").Value = "0.007727"
$ws.Range("E40").Value = "-1.43%"
$ws.Range("D41").Value = "0.1289"
$ws.Range("E41").Value = "-1.90%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.004137"
$ws.Range("E42").Value = "119.87%"
$ws.Range("B43").Value = "Dexo"
$ws.Range("C43").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D43").Value = "0.007580"
$ws.Range("E43").Value = "16.59%"
$ws.Range("D44").Value = "0.007177"
$ws.Range("E44").Value = "-12.96%"
$ws.Range("D45").Value = "0.3202"
$ws.Range("E45").Value = "-4.33%"
$ws.Range("D46").Value = "0.00006836"
$ws.Range("E46").Value = "1.55%"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").Value = "0.48%"
$ws.Range("D48").Value = "0.2517"
$ws.Range("E48").Value = "-32.31%"
$ws.Range("D49").Value = "0.004229"
$ws.Range("E49").Value = "36.29%"
$ws.Range("D50").Value = "0.00002112"
$ws.Range("E50").Value = "0.48%"
$ws.Range("D51").Value = "0.0002011"
$ws.Range("E51").Value = "0.48%"

# Step 3: Restore the default "Normal" style on those cells so the only
# change persisted to the workbook is the cell value/text, matching the
# original (unstyled) formatting of these data cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
